$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, shifting existing rows 194-204 down to 195-205.
$ws.Rows("194:194").Insert()

# Populate the newly inserted row 194 with the new weekly record.
$ws.Range("A194").Value = 11
$ws.Range("B194").Value = "Vega Monumental Concepción"
$ws.Range("C194").Value = "Bíobío"
$ws.Range("D194").Value = 45267
$ws.Range("E194").Value = 8
$ws.Range("F194").Value = 100112001
$ws.Range("G194").Value = "Berenjena"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 50
$ws.Range("K194").Value = 10000
$ws.Range("L194").Value = 10000
$ws.Range("M194").Value = 10000
$ws.Range("N194").Value = '$/caja 50 unidades'
$ws.Range("O194").Value = "Región de Arica y Parinacota"
$ws.Range("P194").Value = 200
$ws.Range("Q194").Value = 50
$ws.Range("R194").Value = "Hortaliza"
